$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 18.71679738969934
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 34.86219683192438

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 9.983522426115931
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 26.01411105592328
